# This script reorders the per-row values in columns D, J, K, L, M, P
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Precio $/Kg) across data rows 2-19, as described by the commit
# "Fruta / hortaliza, semanal". The other columns (A, B, C, E-I, N, O, Q, R)
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values currently sitting in the
# source row before the edit should end up in the destination row).
$map = @{
    2  = 18
    3  = 16
    4  = 9
    5  = 15
    6  = 13
    7  = 4
    8  = 19
    9  = 5
    10 = 8
    11 = 3
    12 = 10
    13 = 11
    14 = 6
    15 = 2
    16 = 7
    17 = 12
    18 = 17
    19 = 14
}

$cols = @(4, 10, 11, 12, 13, 16)  # D, J, K, L, M, P

# Snapshot the original values for the columns we are going to permute,
# before any writes happen, so the permutation is computed consistently.
$orig = @{}
foreach ($row in 2..19) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $orig[$row] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value2 = $orig[$srcRow][$col]
    }
}
